$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '61.401.79'
$ws.Range("E2").Value = '  +0.77%  '
$ws.Range("D3").Value = '3.383.73'
$ws.Range("E3").Value = '  +2.17%  '
$ws.Range("E4").Value = '  -0.06%  '
$ws.Range("D5").Value = '573.06'
$ws.Range("E5").Value = '  +0.88%  '
$ws.Range("D6").Value = '138.03'
$ws.Range("E6").Value = '  +9.24%  '
$ws.Range("E7").Value = '  +0.03%  '
$ws.Range("D8").Value = '3.384.14'
$ws.Range("E8").Value = '  +2.22%  '
$ws.Range("E9").Value = '  +0.03%  '
$ws.Range("E10").Value = '  +5.08%  '
$ws.Range("E11").Value = '  +5.94%  '
$ws.Range("D12").Value = '0.393'
$ws.Range("E12").Value = '  +5.53%  '
$ws.Range("D13").Value = '3.960.95'
$ws.Range("E13").Value = '  +2.00%  '
$ws.Range("E14").Value = '  +2.48%  '
$ws.Range("E15").Value = '  +4.62%  '
$ws.Range("D16").Value = '3.380.94'
$ws.Range("E16").Value = '  +1.80%  '
$ws.Range("D17").Value = '25.19'
$ws.Range("E17").Value = '  +3.78%  '
$ws.Range("D18").Value = '61.480.49'
$ws.Range("E18").Value = '  +0.69%  '
$ws.Range("D19").Value = '13.95'
$ws.Range("E19").Value = '  +6.21%  '
$ws.Range("B20").Value = 'Polkadot'
$ws.Range("C20").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D20").Value = '5.83'
$ws.Range("E20").Value = '  +4.94%  '
$ws.Range("B21").Value = 'Uniswap'
$ws.Range("C21").Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range("D21").Value = '9.48'
$ws.Range("E21").Value = '  +4.88%  '
$ws.Range("D22").Value = '382.18'
$ws.Range("E22").Value = '  +9.63%  '
$ws.Range("E23").Value = '  +4.68%  '
$ws.Range("D24").Value = '3.518.50'
$ws.Range("E24").Value = '  +2.12%  '
$ws.Range("E25").Value = '  +0.26%  '
$ws.Range("D26").Value = '70.92'
$ws.Range("E26").Value = '  +1.59%  '
$ws.Range("E27").Value = '  +13.61%  '
$ws.Range("D28").Value = '1.67'
$ws.Range("E28").Value = '  +16.30%  '
$ws.Range("E29").Value = '  +10.56%  '
$ws.Range("D30").Value = '1.00'
$ws.Range("E30").Value = '  -0.27%  '
$ws.Range("D31").Value = '8.18'
$ws.Range("E31").Value = '  +4.86%  '
$ws.Range("E32").Value = '  +7.19%  '
$ws.Range("E33").Value = '  +2.13%  '
$ws.Range("E34").Value = '  -0.03%  '
$ws.Range("D35").Value = '3.412.86'
$ws.Range("E35").Value = '  +2.08%  '
$ws.Range("D36").Value = '23.52'
$ws.Range("E36").Value = '  +5.53%  '
$ws.Range("E37").Value = '  +4.33%  '
$ws.Range("E38").Value = '  +5.50%  '
$ws.Range("E39").Value = '  +6.42%  '
$ws.Range("D40").Value = '162.86'
$ws.Range("E40").Value = '  +0.19%  '
$ws.Range("D41").Value = '0.0807'
$ws.Range("E41").Value = '  +7.80%  '
$ws.Range("E42").Value = '  -0.13%  '
$ws.Range("E43").Value = '  +5.70%  '
$ws.Range("D44").Value = '41.52'
$ws.Range("E44").Value = '  +1.64%  '
$ws.Range("B45").Value = 'Mantle'
$ws.Range("C45").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("D45").Value = '0.762'
$ws.Range("E45").Value = '  +2.63%  '
$ws.Range("B46").Value = 'ONDO'
$ws.Range("C46").Value = 'https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo'
$ws.Range("D46").Value = '1.20'
$ws.Range("E46").Value = '  +9.52%  '
$ws.Range("D47").Value = '1.69'
$ws.Range("E47").Value = '  +10.06%  '
$ws.Range("D48").Value = '23.18'
$ws.Range("E48").Value = '  +3.38%  '
$ws.Range("E49").Value = '  +5.24%  '
$ws.Range("D50").Value = '23.24'
$ws.Range("E50").Value = '  +13.17%  '
$ws.Range("D51").Value = '2.41'
$ws.Range("E51").Value = '  +13.24%  '
